$d = $word.ActiveDocument

$lines = @(
    "Добавить светодиоды минут внутрь.",
    "Добавить светодиоды с монтажом в отверстие.",
    "Расположить всё так, чтобы было возможно обрезать плату."
)

foreach ($line in $lines) {
    $p = $d.Paragraphs.Last
    $r = $p.Range
    $r.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $line
}

# The document's "_GoBack" bookmark (last-edit marker) must end up collapsed
# right after the text of the final inserted paragraph, not wrapping it.
# Word COM's Bookmarks.Add snaps a range collapsed exactly at the end of a
# paragraph's last run to the whole run, so we briefly append a marker
# character, anchor the bookmark just before it, then remove the marker.
$lastPara = $d.Paragraphs.Last
$markerLen = 1
$lastPara.Range.InsertAfter("#")

$bmRange = $lastPara.Range.Duplicate
$bmRange.MoveEnd(1, -1)
$bmRange.MoveStart(1, $bmRange.Text.Length - $markerLen)
$bmRange.Collapse(1)

try {
    $existing = $d.Bookmarks.Item("_GoBack")
    $existing.Delete()
} catch {
}
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $lastPara.Range.Duplicate
$markerRange.MoveEnd(1, -1)
$markerRange.MoveStart(1, $markerRange.Text.Length - $markerLen)
$markerRange.Text = ""
